$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two product-name cells (renamed test data)
$ws.Range("A2").Value = "Permen"
$ws.Range("H2").Value = "Edited Candy"

# Update the sheet's saved selection/scroll state: select L2 (this also
# clears the previous topLeftCell="D1" scroll position)
$ws.Range("L2").Select() | Out-Null
